$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column K ("2022") values, keyed by row number.
$values = @{
    4  = 2022
    5  = 24.2
    7  = 25.5
    8  = 22.3
    10 = 18
    11 = 18.899999999999999
    12 = 21.3
    13 = 30.2
    14 = 31.7
    15 = 34.1
    16 = 25.8
    17 = 20
    18 = 12.1
    19 = 10.3
    20 = 15.1
    21 = 12.1
    23 = 25.9
    24 = 23.2
    26 = 25.9
    27 = 48.3
    28 = 24.3
    29 = 28.1
    30 = 25.8
    31 = 27.1
    32 = 20.7
    33 = 24.3
    34 = 19.399999999999999
    35 = 7.5
    36 = 11.4
    37 = 36.5
    38 = 17.8
    39 = 20.3
    40 = 20.5
    41 = 32.200000000000003
    42 = 23.2
    43 = 23.8
    44 = 21
    45 = 18
    46 = 3.2
}

# Category header rows that only need the formatting carried across (no value).
$emptyStyleRows = @(6, 9, 22, 25)

for ($r = 4; $r -le 46; $r++) {
    $src = $ws.Range("I$r")
    $dst = $ws.Range("K$r")
    $src.Copy($dst) | Out-Null

    if ($values.ContainsKey($r)) {
        $dst.Value = $values[$r]
    }
}

# Update the saved selection to match the post-edit state.
$ws.Range("L12").Select() | Out-Null
